$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores plain numeric-looking text (e.g. "308.12",
# "1.802.48"). The COM layer auto-detects plain-number-shaped strings and
# coerces them to real numbers (stripping significant trailing zeros, e.g.
# "5.260" becomes 5.26). Forcing NumberFormat to Text ("@") before the
# assignment keeps every updated price cell as literal text, matching the
# original inlineStr cells. The style is reset back to Normal afterwards so
# no stray number-format override is left behind on the cell.
$priceCells = @("D2", "D3", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.451.79"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").Value = "1.803.59"
$ws.Range("E3").Value = "  -2.55%  "
$ws.Range("E4").Value = "  +0.80%  "
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "308.12"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("D7").Value = "0.4558"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("D8").Value = "0.3657"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").Value = "0.07115"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").Value = "0.8758"
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("D11").Value = "0.07722"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").Value = "1.809.08"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "5.260"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "6.347"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").Value = "85.88"
$ws.Range("E16").Value = "  -5.71%  "
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "0.000008566"
$ws.Range("E18").Value = "  -3.86%  "
$ws.Range("D19").Value = "1.008"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "26.503.92"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("D21").Value = "14.23"
$ws.Range("E21").Value = "  -2.93%  "
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").Value = "1.988"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").Value = "150.72"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("D27").Value = "2.007"
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("D28").Value = "112.39"
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").Value = "4.813"
$ws.Range("E29").Value = "  -4.76%  "
$ws.Range("D30").Value = "0.08650"
$ws.Range("E30").Value = "  -1.77%  "
$ws.Range("D31").Value = "3.038"
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").Value = "0.7270"
$ws.Range("E32").Value = "  -4.75%  "
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("E34").Value = "  -4.81%  "
$ws.Range("D35").Value = "1.007"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").Value = "2.530"
$ws.Range("E36").Value = "  -6.95%  "
$ws.Range("D37").Value = "1.077"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "0.01927"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("D40").Value = "2.882"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").Value = "6.913"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Value = "0.4973"
$ws.Range("E42").Value = "  -2.34%  "
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("D46").Value = "0.4584"
$ws.Range("E46").Value = "  -4.09%  "
$ws.Range("D47").Value = "101.49"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").Value = "9.901"
$ws.Range("E48").Value = "  -3.59%  "
$ws.Range("D49").Value = "1.586"
$ws.Range("E49").Value = "  -2.87%  "
$ws.Range("D50").Value = "0.05982"
$ws.Range("D51").Value = "63.77"
$ws.Range("E51").Value = "  -2.28%  "

# Restore default styling on the cells we touched above so no stray
# number-format override is left on the cell (matches the source file,
# where these cells carry no explicit style).
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
